# Generate Report for Handback
# - Status for the a49575e6.. (zh-cn / de-de) handback row flips from
#   "Ready for handoff" to "Handback transform failed" (this status text is a
#   single shared string used by the Overview sheet as well as both locale
#   sheets, so a global replace keeps every reference in sync).
# - The Error Detail (col P) cell for that same row gets a failure message
#   recorded for each locale.
# - Column P is widened to fit the new message text.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Flip the status text everywhere it is used (Overview!E3/F3, zh-cn!C3, de-de!C3)
$ws1.Cells.Replace("Ready for handoff", "Handback transform failed", 1)
$ws2.Cells.Replace("Ready for handoff", "Handback transform failed", 1)
$ws3.Cells.Replace("Ready for handoff", "Handback transform failed", 1)

# Record the handback/handoff file-name mismatch error per locale
$ws2.Range("P3").Value = "Handback file name: gsiasq4g.k3r is different with handoff file name: a49575e6-5d85-402b-a908-e03cf43dbd03.03951f21259fc8d25f03fdc4016274765160b3bd.zh-cn."
$ws3.Range("P3").Value = "Handback file name: gsiasq4g.k3r is different with handoff file name: a49575e6-5d85-402b-a908-e03cf43dbd03.03951f21259fc8d25f03fdc4016274765160b3bd.de-de."

# Widen the Error Detail column (P) so the new message is readable
$ws2.Columns.Item(16).ColumnWidth = 39.17
$ws3.Columns.Item(16).ColumnWidth = 39.17
